# RPS Qualifying Source Definitions - update for South Korea data source
# (per commit "updated south korea files")

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "About": replace the old US-centric source note with the new
# South Korea (Korea Energy Agency / RPS) source block, and shrink the
# long explanatory "Notes" paragraph down to the two short KR notes.
# -------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Remove the old long-form note rows (rows 9-18): the "Each U.S. state..."
# paragraph plus the BAU explanation bullets. This leaves row 8 ("Notes")
# directly under row 5.
$wsAbout.Rows("9:18").Delete()

# Insert one blank row above the "Notes" header so it lands on row 9,
# leaving two free rows (6 and 7) between the source block and it -
# matching the target layout.
$wsAbout.Rows("8").Insert()

# Fill in the new "Source" block. Cells are written in this particular
# order so that new shared-string entries are appended in the same
# sequence the original authors entered them in.
$wsAbout.Range("B7").Value = "https://www.knrec.or.kr/business/rps_guide.aspx"
$wsAbout.Range("B4").Value = "RPS"
$wsAbout.Range("B5").Value = "Korea Energy Agency"
$wsAbout.Range("B5").Style = "Hyperlink"
$wsAbout.Range("B6").Value = "Renewable energy certificates"
$wsAbout.Range("A11").Value = "There is no difference by region in South Korea"
$wsAbout.Range("A10").Value = "We allocated boolean values based on the list provided by Korea Energy Agency."

# Restore the gray "source" highlight + bold look on the value cell next
# to "Source:" (kept from the original template, just re-applied here).
$wsAbout.Range("B4").Font.Bold = $true
$wsAbout.Range("B4").Interior.ThemeColor = 2
$wsAbout.Range("B4").Interior.TintAndShade = -0.249977111117893

# -------------------------------------------------------------------------
# Sheet "RQSD-BRQSD" (BAU qualifying sources): nuclear no longer qualifies.
# -------------------------------------------------------------------------
$wsBau = $wb.Worksheets.Item("RQSD-BRQSD")
$wsBau.Range("B4").Value = 0

# -------------------------------------------------------------------------
# Sheet "RQSD-RQSD" (non-BAU qualifying sources): hydro and biomass now
# qualify, per the Korea Energy Agency RPS list.
# -------------------------------------------------------------------------
$wsRqsd = $wb.Worksheets.Item("RQSD-RQSD")
$wsRqsd.Range("B5").Value = 1
$wsRqsd.Range("B9").Value = 1

# -------------------------------------------------------------------------
# Restore the selection state on each sheet (matches the saved workbook),
# finishing on "About" so it stays the active tab.
# -------------------------------------------------------------------------
$wsBau.Range("D5").Select()
$wsRqsd.Range("B4").Select()
$wsAbout.Activate()
$wsAbout.Range("A11").Select()
